# Updated cryptos list on Wed Oct 18 16:42:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to be treated/stored as text so Excel does not
    # auto-convert number-like strings (e.g. "211.93") into real numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.329.44"
Set-TextValue $ws.Range("E2") "  -0.85%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.570.23"
Set-TextValue $ws.Range("E3") "  -0.64%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "211.93"
Set-TextValue $ws.Range("E5") "  -0.24%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.489"
Set-TextValue $ws.Range("E6") "  -0.52%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  +0.03%  "

# Row 8 - OKB
Set-TextValue $ws.Range("D8") "44.55"
Set-TextValue $ws.Range("E8") "  -4.80%  "

# Row 9 - Solana
Set-TextValue $ws.Range("D9") "23.65"
Set-TextValue $ws.Range("E9") "  -2.01%  "

# Row 10 - Cardano
Set-TextValue $ws.Range("D10") "0.245"
Set-TextValue $ws.Range("E10") "  -0.98%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0586"
Set-TextValue $ws.Range("E11") "  -0.99%  "

# Row 12 - TRON
Set-TextValue $ws.Range("D12") "0.0895"
Set-TextValue $ws.Range("E12") "  +1.50%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "1.797.45"
Set-TextValue $ws.Range("E13") "  -0.44%  "

# Row 14 - WrappedEther
Set-TextValue $ws.Range("D14") "1.573.99"
Set-TextValue $ws.Range("E14") "  -0.46%  "

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") "3.68"
Set-TextValue $ws.Range("E15") "  -0.63%  "

# Row 16 - now WrappedBTC (was Polygon)
Set-TextValue $ws.Range("B16") "WrappedBTC"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "28.356.78"
Set-TextValue $ws.Range("E16") "  -0.69%  "

# Row 17 - now Polygon (was WrappedBTC)
Set-TextValue $ws.Range("B17") "Polygon"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D17") "0.515"
Set-TextValue $ws.Range("E17") "  -1.67%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "61.46"
Set-TextValue $ws.Range("E18") "  -1.49%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "229.76"
Set-TextValue $ws.Range("E19") "  +0.27%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.40"
Set-TextValue $ws.Range("E20") "  -0.09%  "

# Row 21 - ShibaInu
Set-TextValue $ws.Range("D21") "0.0₃0683"
Set-TextValue $ws.Range("E21") "  -1.58%  "

# Row 22 - Dai
Set-TextValue $ws.Range("E22") "  -0.03%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "3.95"
Set-TextValue $ws.Range("E23") "  +0.66%  "

# Row 24 - Avalanche
Set-TextValue $ws.Range("D24") "9.02"
Set-TextValue $ws.Range("E24") "  -1.56%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("E25") "  +0.40%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "151.02"
Set-TextValue $ws.Range("E26") "  -0.08%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "14.88"
Set-TextValue $ws.Range("E27") "  -0.88%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "6.35"
Set-TextValue $ws.Range("E28") "  -1.77%  "

# Row 29 - Stellar
Set-TextValue $ws.Range("D29") "0.103"
Set-TextValue $ws.Range("E29") "  -1.24%  "

# Row 30 - BinanceUSD
Set-TextValue $ws.Range("E30") "  +0.00%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("E31") "  +3.27%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("E32") "  -3.75%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("E33") "  -1.01%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "3.11"
Set-TextValue $ws.Range("E34") "  -0.74%  "

# Row 35 - Maker
Set-TextValue $ws.Range("D35") "1.383.23"
Set-TextValue $ws.Range("E35") "  -1.11%  "

# Row 36 - TrustWalletToken
Set-TextValue $ws.Range("E36") "  +5.27%  "

# Row 37 - LidoDAOToken
Set-TextValue $ws.Range("D37") "1.50"
Set-TextValue $ws.Range("E37") "  -3.16%  "

# Row 38 - HuobiToken
Set-TextValue $ws.Range("E38") "  -0.15%  "

# Row 39 - MXToken
Set-TextValue $ws.Range("E39") "  +1.54%  "

# Row 40 - VeChain
Set-TextValue $ws.Range("E40") "  -1.88%  "

# Row 41 - ImmutableX
Set-TextValue $ws.Range("D41") "0.518"
Set-TextValue $ws.Range("E41") "  -2.63%  "

# Row 42 - PaxDollar
Set-TextValue $ws.Range("E42") "  -0.12%  "

# Row 43 - RenderToken
Set-TextValue $ws.Range("D43") "1.89"
Set-TextValue $ws.Range("E43") "  +1.57%  "

# Row 44 - ARBITRUM
Set-TextValue $ws.Range("D44") "0.784"
Set-TextValue $ws.Range("E44") "  -1.58%  "

# Row 45 - Kaspa
Set-TextValue $ws.Range("E45") "  -0.04%  "

# Row 46 - FraxShare
Set-TextValue $ws.Range("D46") "5.37"
Set-TextValue $ws.Range("E46") "  -4.22%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "62.16"
Set-TextValue $ws.Range("E47") "  -1.29%  "

# Row 48 - WEMIXToken
Set-TextValue $ws.Range("D48") "0.920"
Set-TextValue $ws.Range("E48") "  -5.92%  "

# Row 49 - RocketPoolETH
Set-TextValue $ws.Range("D49") "1.709.86"
Set-TextValue $ws.Range("E49") "  -0.37%  "

# Row 50 - mCoin
Set-TextValue $ws.Range("E50") "  +0.94%  "

# Row 51 - Quant
Set-TextValue $ws.Range("D51") "85.14"
Set-TextValue $ws.Range("E51") "  -1.02%  "
